$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Deal sheet: move selection from L5 -> C2
# ------------------------------------------------------------------
$dealWs = $wb.Worksheets.Item("Deal")
$dealWs.Activate()
$dealWs.Range("C2").Select()

# ------------------------------------------------------------------
# 2) CustomSDG sheet: widen selection from A1:B1 -> A1:D1
# ------------------------------------------------------------------
$sdgWs = $wb.Worksheets.Item("CustomSDG")
$sdgWs.Activate()
$sdgWs.Range("A1:D1").Select()

# ------------------------------------------------------------------
# 3) Insert new "DealRequestTracker" sheet right before "ToggleBtn"
# ------------------------------------------------------------------
$toggleBtn = $wb.Worksheets.Item("ToggleBtn")
$tracker = $wb.Worksheets.Add($toggleBtn)
$tracker.Name = "DealRequestTracker"

# Values are entered in the exact sequence needed to reproduce the
# original shared-strings ordering (first-use order) from the diff.
$tracker.Range("A1").Value = "Variable_Name"
$tracker.Range("A2").Value = "OPENQA1"
$tracker.Range("A3").Value = "CLOSEDQA1"
$tracker.Range("B1").Value = "Request_Tracker_ID"
$tracker.Range("D1").Value = "Request"
$tracker.Range("E1").Value = "Status"
$tracker.Range("E2").Value = "Open"
$tracker.Range("B2").Value = "RT"
$tracker.Range("C1").Value = "Date_Requested"
$tracker.Range("D3").Value = "Finance Related"
$tracker.Range("C2").Value = "'12/22/2020"
$tracker.Range("C2").NumberFormat = "mm-dd-yy"
$tracker.Range("C2").Copy() | Out-Null
$tracker.Range("C3").PasteSpecial(-4122) | Out-Null
$tracker.Range("C3").Value = "'12/7/2020"
$tracker.Range("D2").Value = "IT Related"
$tracker.Range("B3").Value = "RT"

# Column widths to roughly match the authored layout
$tracker.Columns("A").ColumnWidth = 14.85546875
$tracker.Columns("B").ColumnWidth = 18.85546875
$tracker.Columns("C:D").ColumnWidth = 15.28515625

$tracker.Range("D3").Select()

# Make the new sheet the active / visible tab (mirrors tabSelected moving here)
$tracker.Activate()

Write-Host "DealRequestTracker sheet created"
